$wb = $excel.ActiveWorkbook

# 1. Rename sheet "fond privat" -> "fond" and fix the title cells that
#    used to spell out the old sheet name.
$wsFond = $wb.Worksheets.Item("fond privat")
$wsFond.Cells.Item(1, 1).Value = "fond"
$wsFond.Cells.Item(2, 2).Value = "fond"
$wsFond.Name = "fond"

# 2. Resize the first two columns on the "fond" sheet to fit the shorter text.
$wsFond.Columns.Item(1).ColumnWidth = 10.569
$wsFond.Columns.Item(2).ColumnWidth = 6.998

# 3. Rename "...Using previous/current value:..." to "...amount:..." everywhere.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $null = $used.Replace("Using previous value:", "Using previous amount:")
    $null = $used.Replace("Using current value:", "Using current amount:")
}
